$wb = $excel.ActiveWorkbook

# Helper to add a sheet at the end
function Add-SheetAtEnd($name) {
    $count = $wb.Worksheets.Count
    $ws = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets($count))
    $ws.Name = $name
    return $ws
}

# --- Sheet: 05-36 Pakete über das Terminal ---
$ws12 = Add-SheetAtEnd "05-36 Pakete über das Terminal"
$ws12.Range("A1").Value = "Frage"
$ws12.Range("B1").Value = "Antwort"
$ws12.Range("A2").Value = "UBUNTU: Update die Paketlisten"
$ws12.Range("B2").Value = "apt update"
$ws12.Range("B3").Value = "apt upgrade"
$ws12.Range("B4").Value = "apt dist-upgrade"
$ws12.Range("A3").Value = "Aktualisiere und installiere die neueste Version aller Pakete ohne neue Abhängigkeiten"
$ws12.Range("A4").Value = "Aktualisiere und installiere die neueste Version aller Pakete mit neuen Abhängigkeiten"
$ws12.Range("A5").Value = "UBUNTU: Installiere htop"
$ws12.Range("B5").Value = "apt install htop"
$ws12.Columns.Item(1).ColumnWidth = 79.7369791667
$ws12.Columns.Item(2).ColumnWidth = 14.7369791667
$ws12.Range("A6").Select() | Out-Null

# --- Sheet: 05-37 Pakete + Quellen suchen ---
$ws13 = Add-SheetAtEnd "05-37 Pakete + Quellen suchen"
$ws13.Range("A1").Value = "Frage"
$ws13.Range("B1").Value = "Antwort"
$ws13.Range("B2").Value = "apt-get install php"
$ws13.Range("A3").Value = "Alternativer Debianbefehl (nicht apt, nicht apt-get) um php zu installieren"
$ws13.Range("A2").Value = "Alternativer Debianbefehl (nicht apt, nicht aptitude) um php zu installieren"
$ws13.Range("B3").Value = "aptitude install php"
$ws13.Range("A4").Value = "Suche in den Paketlisten nach php"
$ws13.Range("B4").Value = "apt search php"
$ws13.Range("A5").Value = "Website für Ubuntu packages"
$ws13.Range("B5").Value = "packages.ubuntu.com"
$ws13.Range("A6").Value = "Website für Debian packages"
$ws13.Range("B6").Value = "packages.debian.org"
$ws13.Range("A7").Value = "Wo liegen die Konfigurationsdatein für den Paketmanager apt?"
$ws13.Range("B7").Value = "/etc/apt/sources.list"
$ws13.Range("A8").Value = "Wo liegen die Konfigurationsdateien für selbst installierte Pakete mit dpkg?"
$ws13.Range("B8").Value = "/etc/apt/sources.list.d/"
$ws13.Columns.Item(1).ColumnWidth = 68.1666666667
$ws13.Columns.Item(2).ColumnWidth = 16.4518229167
$ws13.Range("A9").Select() | Out-Null

# --- Sheet: 05-38 Skype nachinstallieren ---
$ws14 = Add-SheetAtEnd "05-38 Skype nachinstallieren"
$ws14.Range("A1").Value = "Frage"
$ws14.Range("B1").Value = "Antwort"
$ws14.Range("A2").Select() | Out-Null
